# Weekly refresh of the "Ciruela" (plum) price sheet: two new price
# observations (row 90 and 91, dated 2022-03-22 / serial 44642) are
# published at the top of the data block, pushing the previously-existing
# rows 90-108 down to 92-110 (dimension grows from A1:T108 to A1:T110).
# All of the shifted rows keep their original values untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 90; this shifts the
# existing rows 90-108 down to 92-110, exactly like Excel's own
# Rows.Insert (shift-down) behaviour.
$ws.Rows.Item(90).Insert()
$ws.Rows.Item(90).Insert()

# New rows to populate, in column A..T order.
$newRows = @(
    @{ Row = 90; Values = @(5, "Macroferia Regional de Talca", "Maule", 44642, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Primera", 230, 8000, 8000, 8000, "$/bandeja 18 kilos granel", "Región del Maule", 444, 18) },
    @{ Row = 91; Values = @(5, "Macroferia Regional de Talca", "Maule", 44642, 7, "Fruta", 100103, "Frutos de hueso (carozo)", 100103002, "Ciruela", "Angeleno", "Segunda", 200, 6000, 6000, 6000, "$/bandeja 18 kilos granel", "Región del Maule", 333, 18) }
)

foreach ($entry in $newRows) {
    $r = $entry.Row
    $vals = $entry.Values
    for ($col = 1; $col -le $vals.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $vals[$col - 1]
    }
}
